# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Only column G (header "K") values for rows 2-7 change; update them directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 1
